# Auto-generated edit script applying the Sagittarius_Profits market-data refresh
# Updates currentAveragePrice / Leve price / profit columns (H:N) across all 8 job sheets,
# mirroring a scheduled market-data pull. Values only; no formulas, no style changes.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1178.2858
$ws.Range("I19").Value = 702.6667
$ws.Range("K19").Value = 702.6667
$ws.Range("M19").Value = -527.6667
$ws.Range("H40").Value = 2208.4092
$ws.Range("I40").Value = 1881.1818
$ws.Range("J40").Value = 2535.6365
$ws.Range("K40").Value = 1881.1818
$ws.Range("L40").Value = 2535.6365
$ws.Range("M40").Value = -1706.1818
$ws.Range("N40").Value = -2885.6365
$ws.Range("H70").Value = 1610.875
$ws.Range("I70").Value = 999
$ws.Range("J70").Value = 1698.2858
$ws.Range("K70").Value = 2997
$ws.Range("L70").Value = 5094.857400000001
$ws.Range("M70").Value = -2727
$ws.Range("N70").Value = -5634.857400000001
$ws.Range("H73").Value = 1610.875
$ws.Range("I73").Value = 999
$ws.Range("J73").Value = 1698.2858
$ws.Range("K73").Value = 2997
$ws.Range("L73").Value = 5094.857400000001
$ws.Range("M73").Value = -2061
$ws.Range("N73").Value = -6966.857400000001
$ws.Range("H106").Value = 150000
$ws.Range("I106").Value = 150000
$ws.Range("K106").Value = 150000
$ws.Range("M106").Value = -149369
$ws.Range("H116").Value = 7999.5
$ws.Range("I116").Value = 7999
$ws.Range("K116").Value = 7999
$ws.Range("M116").Value = -4557

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2486499.8
$ws.Range("J32").Value = 779803.9
$ws.Range("L32").Value = 779803.9
$ws.Range("N32").Value = -780377.9
$ws.Range("H38").Value = 63779.57
$ws.Range("I38").Value = 73742.836
$ws.Range("K38").Value = 73742.836
$ws.Range("M38").Value = -73275.836
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("K61").Value = 3000
$ws.Range("M61").Value = -2788
$ws.Range("H74").Value = 760.1667
$ws.Range("I74").Value = 712.2
$ws.Range("K74").Value = 712.2
$ws.Range("M74").Value = 161.8
$ws.Range("H77").Value = 760.1667
$ws.Range("I77").Value = 712.2
$ws.Range("K77").Value = 3561
$ws.Range("M77").Value = 807
$ws.Range("H109").Value = 125000
$ws.Range("J109").Value = 125000
$ws.Range("L109").Value = 125000
$ws.Range("N109").Value = -127774
$ws.Range("H112").Value = 26212.5
$ws.Range("J112").Value = 26212.5
$ws.Range("L112").Value = 26212.5
$ws.Range("N112").Value = -29166.5
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 148
$ws.Range("I22").Value = 160.25
$ws.Range("J22").Value = 99
$ws.Range("K22").Value = 160.25
$ws.Range("L22").Value = 99
$ws.Range("M22").Value = 12.75
$ws.Range("N22").Value = -445
$ws.Range("H62").Value = 39999
$ws.Range("J62").Value = 39999
$ws.Range("L62").Value = 39999
$ws.Range("N62").Value = -41371
$ws.Range("H65").Value = 39999
$ws.Range("J65").Value = 39999
$ws.Range("L65").Value = 119997
$ws.Range("N65").Value = -126861
$ws.Range("H105").Value = 2245.4167
$ws.Range("J105").Value = 1862.25
$ws.Range("L105").Value = 1862.25
$ws.Range("N105").Value = -5356.25
$ws.Range("H134").Value = 2677.4443
$ws.Range("I134").Value = 2677.4443
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8032.3329
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5497.3329
$ws.Range("N134").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1161.4445
$ws.Range("J5").Value = 762
$ws.Range("L5").Value = 762
$ws.Range("N5").Value = -986
$ws.Range("H8").Value = 2820
$ws.Range("I8").Value = 2820
$ws.Range("K8").Value = 2820
$ws.Range("M8").Value = -2680
$ws.Range("H12").Value = 3361.6667
$ws.Range("I12").Value = 2542.5
$ws.Range("K12").Value = 2542.5
$ws.Range("M12").Value = -2372.5
$ws.Range("H25").Value = 700
$ws.Range("I25").Value = 700
$ws.Range("K25").Value = 700
$ws.Range("M25").Value = -526
$ws.Range("H31").Value = 998.3333
$ws.Range("I31").Value = 999
$ws.Range("J31").Value = 997
$ws.Range("K31").Value = 999
$ws.Range("L31").Value = 997
$ws.Range("M31").Value = -704
$ws.Range("N31").Value = -1587
$ws.Range("H34").Value = 998.3333
$ws.Range("I34").Value = 999
$ws.Range("J34").Value = 997
$ws.Range("K34").Value = 999
$ws.Range("L34").Value = 997
$ws.Range("M34").Value = -797
$ws.Range("N34").Value = -1401
$ws.Range("H58").Value = 2001.5714
$ws.Range("I58").Value = 2001.5714
$ws.Range("K58").Value = 2001.5714
$ws.Range("M58").Value = -1798.5714
$ws.Range("H62").Value = 6159.6
$ws.Range("J62").Value = 6159.6
$ws.Range("L62").Value = 6159.6
$ws.Range("N62").Value = -7407.6
$ws.Range("H65").Value = 6159.6
$ws.Range("J65").Value = 6159.6
$ws.Range("L65").Value = 30798
$ws.Range("N65").Value = -37038
$ws.Range("H86").Value = 19998.5
$ws.Range("I86").Value = 19998
$ws.Range("K86").Value = 19998
$ws.Range("M86").Value = -18875
$ws.Range("H89").Value = 19998.5
$ws.Range("I89").Value = 19998
$ws.Range("K89").Value = 99990
$ws.Range("M89").Value = -94374
$ws.Range("H134").Value = 3033.3333
$ws.Range("I134").Value = 3033.3333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9099.999899999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6564.999899999999
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2001.5714
$ws.Range("I136").Value = 2001.5714
$ws.Range("K136").Value = 6004.7142
$ws.Range("M136").Value = -3454.7142

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86.42856999999999
$ws.Range("J2").Value = 99
$ws.Range("L2").Value = 594
$ws.Range("N2").Value = -820
$ws.Range("H5").Value = 1124.75
$ws.Range("I5").Value = 333
$ws.Range("K5").Value = 999
$ws.Range("M5").Value = -887
$ws.Range("H114").Value = 2304.375
$ws.Range("I114").Value = 1336.6666
$ws.Range("J114").Value = 2885
$ws.Range("K114").Value = 4009.9998
$ws.Range("L114").Value = 8655
$ws.Range("M114").Value = -755.9998000000001
$ws.Range("N114").Value = -15163
$ws.Range("H118").Value = 2400
$ws.Range("I118").Value = 2400
$ws.Range("K118").Value = 7200
$ws.Range("M118").Value = -5957
$ws.Range("H135").Value = 1124.75
$ws.Range("I135").Value = 333
$ws.Range("K135").Value = 2997
$ws.Range("M135").Value = -462
$ws.Range("H137").Value = 1819.6666
$ws.Range("I137").Value = 1068.6
$ws.Range("K137").Value = 3205.8
$ws.Range("M137").Value = 1894.2
$ws.Range("H140").Value = 14859.833
$ws.Range("J140").Value = 16899.5
$ws.Range("L140").Value = 50698.5
$ws.Range("N140").Value = -61058.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -830
$ws.Range("N9").ClearContents()
$ws.Range("H10").Value = 13125
$ws.Range("I10").Value = 11000
$ws.Range("K10").Value = 11000
$ws.Range("M10").Value = -10831
$ws.Range("H113").Value = 655.25
$ws.Range("I113").Value = 249
$ws.Range("J113").Value = 899
$ws.Range("K113").Value = 249
$ws.Range("L113").Value = 899
$ws.Range("M113").Value = 1921
$ws.Range("N113").Value = -5239
$ws.Range("H121").Value = 150000
$ws.Range("J121").Value = 150000
$ws.Range("L121").Value = 150000
$ws.Range("N121").Value = -153494

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 58999
$ws.Range("J64").Value = 58999
$ws.Range("L64").Value = 58999
$ws.Range("N64").Value = -59495
$ws.Range("H67").Value = 58999
$ws.Range("J67").Value = 58999
$ws.Range("L67").Value = 58999
$ws.Range("N67").Value = -60715
$ws.Range("H126").Value = 5221.606
$ws.Range("I126").Value = 4935.905
$ws.Range("J126").Value = 5721.5835
$ws.Range("K126").Value = 14807.715
$ws.Range("L126").Value = 17164.7505
$ws.Range("M126").Value = -12337.715
$ws.Range("N126").Value = -22104.7505
$ws.Range("H136").Value = 3321.2104
$ws.Range("I136").Value = 3621.9285
$ws.Range("K136").Value = 10865.7855
$ws.Range("M136").Value = -8315.7855
